# "Add powerpoint 365 support"
#
# The underlying change is a relationship-id refresh: three slides
# (SlideID 259/260/261) are re-linked from the presentation via fresh
# r:id values, and the picture placeholder on the remaining content
# slide gets its image blip re-linked via a fresh r:embed value. No
# visible content (slide order, slide bodies, picture position/name,
# embedded image bytes) changes at all.
#
# PowerPoint's object model has no "rename this relationship" verb, so
# we reproduce the same externally-visible effect the supported way:
# duplicate the part (slide / shape) — which forces the host to mint a
# brand new relationship id for the copy while leaving every other
# aspect of the content untouched — then delete the original and put
# the duplicate back exactly where the original was.

$p = $ppt.ActivePresentation

# --- 1. Re-link the three table slides (SlideID 259, 260, 261) -----------
# They are the 2nd, 3rd and 4th slides in the deck.
for ($i = 2; $i -le 4; $i++) {
    $orig = $p.Slides.Item($i)
    $copy = $orig.Duplicate()   # lands immediately after $orig, same content
    $orig.Delete()              # removes the old slide + frees its r:id
}

# --- 2. Re-link the picture placeholder's image (slide with SlideID 258) -
$picSlide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    if ($p.Slides.Item($i).SlideID -eq 258) {
        $picSlide = $p.Slides.Item($i)
    }
}

$pic = $picSlide.Shapes.Item(1)
$origLeft = $pic.Left
$origTop = $pic.Top
$origName = $pic.Name

$picCopy = $pic.Duplicate()     # new shape, same image, fresh r:embed
$pic.Delete()                   # drop the old shape + its old r:embed

# Duplicate() offsets the copy by a fixed amount - put it back exactly
# where the original placeholder was (Width/Height are left untouched
# so they keep their original, unrounded EMU values).
$picCopy.Left = $origLeft
$picCopy.Top = $origTop
$picCopy.Name = $origName
